$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..53, columns A..E (date serial, y0, y0_forecast, y1, y1_forecast)
$data = @(
    @(39400, 2007, 0.4235526809466261, 2008, $null),
    @(39583, 2008, $null, 2009, $null),
    @(39765, 2008, -0.571807692896309, 2009, $null),
    @(39948, 2009, $null, 2010, $null),
    @(40130, 2009, 0.3486139762224783, 2010, $null),
    @(40310, 2010, -1.0538757111625, 2011, -0.03775702182785556),
    @(40494, 2010, -0.1384957661262898, 2011, 0.592211799485276),
    @(40676, 2011, 1.570920254756558, 2012, 0.485781149991249),
    @(40862, 2011, 1.566479473280147, 2012, 0.5295895589954247),
    @(41044, 2012, 0.878183952843048, 2013, 0.6274179984581574),
    @(41228, 2012, 0.7307568962936939, 2013, 0.5955791956549161),
    @(41409, 2013, 0.4413109953541605, 2014, 0.6943226143418357),
    @(41592, 2013, 0.818818812164257, 2014, 0.9193568360546411),
    @(41774, 2014, 0.9083436352771646, 2015, 0.8175300924317952),
    @(41957, 2014, 0.9180054319587239, 2015, 0.9697679806505821),
    @(42137, 2015, 1.93704103170067, 2016, 0.8819846436028733),
    @(42321, 2015, 1.984684278296656, 2016, 0.970573649360662),
    @(42503, 2016, 1.58745197360306, 2017, 0.9174072565958813),
    @(42689, 2016, 1.755995812646982, 2017, 1.083941060573212),
    @(42867, 2017, 1.276746817047392, 2018, 0.9840692158344266),
    @(43053, 2017, 1.946965557828384, 2018, 1.343559319682996),
    @(43145, 2018, 0.9903126026122644, 2019, 1.080096276728471),
    @(43235, 2018, 1.116346046342809, 2019, 1.076386909629967),
    @(43326, 2018, 1.219509802703822, 2019, 1.125010763074807),
    @(43418, 2018, 1.06432145354225, 2019, 0.8791375467670726),
    @(43510, 2019, 0.5566732908481198, 2020, 0.9336434107631186),
    @(43600, 2019, 1.437659299153027, 2020, 1.280675138581411),
    @(43691, 2019, 1.315960672561256, 2020, 1.194139179128939),
    @(43783, 2019, 1.361817904277696, 2020, 1.226035857429442),
    @(43875, 2020, 1.08139044988087, 2021, 1.151701057323207),
    @(43966, 2020, -1.741236505435428, 2021, -0.2374583495023508),
    @(44068, 2020, -5.210209911466245, 2021, -0.6076997714017818),
    @(44159, 2020, -4.352425014431304, 2021, -1.242805832929039),
    @(44251, 2021, -2.740233405906423, 2022, 0.509185690804248),
    @(44341, 2021, -3.281345655000223, 2022, -0.1411395765833157),
    @(44432, 2021, -1.897775264882628, 2022, 0.6993852707122938),
    @(44525, 2021, -1.761645650979182, 2022, 3.46909983288044),
    @(44617, 2022, 3.920818017292138, 2023, 1.110570829887347),
    @(44706, 2022, 4.407642498961351, 2023, 1.101211681542535),
    @(44798, 2022, 4.787836378515364, 2023, 1.219524838989372),
    @(44890, 2022, 5.20787683103745, 2023, 1.517423464826884),
    @(44981, 2023, 0.04995351332486297, 2024, 0.861118176709641),
    @(45071, 2023, -1.129337961135779, 2024, 0.5661352122229735),
    @(45163, 2023, -0.49899188013105, 2024, 0.6655557872750961),
    @(45254, 2023, -0.9008525709169546, 2024, 0.2090788898015949),
    @(45345, 2024, 0.4162071511232046, 2025, 0.7845554863595172),
    @(45436, 2024, 0.2034622429862187, 2025, 0.6701250445558804),
    @(45534, 2024, 0.5084754301873051, 2025, 0.7125770617453231),
    @(45618, 2024, 0.2738544794132824, 2025, 0.5362040463673612),
    @(45713, 2025, 0.5222112282702973, 2026, 0.8164751105439905),
    @(45800, 2025, 0.9632017291179906, 2026, 0.8985052439231866),
    @(45891, 2025, 0.9584581489103794, 2026, 0.8209425653958213)
)

$rowIndex = 2
foreach ($rowVals in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowVals[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowVals[1]
    if ($null -eq $rowVals[2]) {
        $ws.Cells.Item($rowIndex, 3).ClearContents()
    } else {
        $ws.Cells.Item($rowIndex, 3).Value = $rowVals[2]
    }
    $ws.Cells.Item($rowIndex, 4).Value = $rowVals[3]
    if ($null -eq $rowVals[4]) {
        $ws.Cells.Item($rowIndex, 5).ClearContents()
    } else {
        $ws.Cells.Item($rowIndex, 5).Value = $rowVals[4]
    }
    $rowIndex++
}

# New row 53 needs column A to carry the same date style as other A cells.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
